$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12419
$ws1.Range("F4").Value = 2020
$ws1.Range("F8").Value = 12357
$ws1.Range("F9").Value = 3015
$ws1.Range("F10").Value = 524
$ws1.Range("F13").Value = 15
$ws1.Range("F17").Value = 6041
$ws1.Range("F18").Value = 141
$ws1.Range("F19").Value = 3591

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12419
$ws4.Range("F4").Value = 2020
$ws4.Range("F9").Value = 12357
$ws4.Range("F10").Value = 3015
$ws4.Range("F11").Value = 524
$ws4.Range("F14").Value = 15
$ws4.Range("F19").Value = 6041
$ws4.Range("F20").Value = 141
$ws4.Range("F21").Value = 3591

$wb.Save()
